$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT, even when the string looks like a plain
# number (e.g. "562.25"), without leaving a lasting NumberFormat/style change
# on the cell (matches the source workbook, where these are inlineStr cells
# with no explicit style).
function Set-TextValue {
    param($cell, [string]$text)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Would otherwise be auto-coerced to a Number by Excel -- force Text,
        # assign, then restore the default "Normal" style so no stray
        # per-cell number format sticks around.
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Range("D2") '62.231.00'
Set-TextValue $ws.Range("E2") '  +1.36%  '
Set-TextValue $ws.Range("D3") '2.420.58'
Set-TextValue $ws.Range("E3") '  +1.77%  '
Set-TextValue $ws.Range("E4") '  -0.10%  '
Set-TextValue $ws.Range("D5") '562.25'
Set-TextValue $ws.Range("E5") '  +1.68%  '
Set-TextValue $ws.Range("D6") '143.80'
Set-TextValue $ws.Range("E6") '  +2.96%  '
Set-TextValue $ws.Range("E7") '  +0.01%  '
Set-TextValue $ws.Range("E8") '  +1.55%  '
Set-TextValue $ws.Range("D9") '2.418.12'
Set-TextValue $ws.Range("E9") '  +1.66%  '
Set-TextValue $ws.Range("D10") '0.109'
Set-TextValue $ws.Range("E10") '  +1.41%  '
Set-TextValue $ws.Range("E11") '  -2.11%  '
Set-TextValue $ws.Range("D12") '5.37'
Set-TextValue $ws.Range("E12") '  +0.26%  '
Set-TextValue $ws.Range("E13") '  +0.37%  '
Set-TextValue $ws.Range("D14") '25.94'
Set-TextValue $ws.Range("E14") '  +1.54%  '
Set-TextValue $ws.Range("E15") '  +2.64%  '
Set-TextValue $ws.Range("D16") '2.858.44'
Set-TextValue $ws.Range("E16") '  +1.77%  '
Set-TextValue $ws.Range("D17") '61.954.42'
Set-TextValue $ws.Range("E17") '  +0.90%  '
Set-TextValue $ws.Range("D18") '2.421.09'
Set-TextValue $ws.Range("E18") '  +1.67%  '
Set-TextValue $ws.Range("E19") '  +3.34%  '
Set-TextValue $ws.Range("E20") '  +0.96%  '
Set-TextValue $ws.Range("D21") '323.86'
Set-TextValue $ws.Range("E21") '  +0.90%  '
Set-TextValue $ws.Range("E22") '  +0.91%  '
Set-TextValue $ws.Range("E23") '  +0.05%  '
Set-TextValue $ws.Range("D24") '65.55'
Set-TextValue $ws.Range("E24") '  +1.98%  '
Set-TextValue $ws.Range("E25") '  -2.16%  '
Set-TextValue $ws.Range("D26") '8.92'
Set-TextValue $ws.Range("E26") '  +0.81%  '
Set-TextValue $ws.Range("D27") '584.91'
Set-TextValue $ws.Range("E27") '  +12.37%  '
Set-TextValue $ws.Range("D28") '2.540.33'
Set-TextValue $ws.Range("E28") '  +1.78%  '
Set-TextValue $ws.Range("D29") '0.997'
Set-TextValue $ws.Range("E29") '  -0.01%  '
Set-TextValue $ws.Range("D30") '0.0₃0944'
Set-TextValue $ws.Range("E30") '  +4.49%  '
Set-TextValue $ws.Range("D31") '1.46'
Set-TextValue $ws.Range("E31") '  +5.09%  '
Set-TextValue $ws.Range("D32") '8.24'
Set-TextValue $ws.Range("E32") '  +0.99%  '
Set-TextValue $ws.Range("E33") '  +1.18%  '
Set-TextValue $ws.Range("E34") '  +2.43%  '
Set-TextValue $ws.Range("E35") '  +1.73%  '
Set-TextValue $ws.Range("E36") '  +3.72%  '
Set-TextValue $ws.Range("E37") '  +0.04%  '
Set-TextValue $ws.Range("D38") '4.79'
Set-TextValue $ws.Range("E38") '  +2.17%  '
Set-TextValue $ws.Range("D39") '0.384'
Set-TextValue $ws.Range("E39") '  +1.60%  '
Set-TextValue $ws.Range("D40") '152.84'
Set-TextValue $ws.Range("E40") '  +4.23%  '
Set-TextValue $ws.Range("D41") '18.68'
Set-TextValue $ws.Range("E41") '  +0.96%  '
Set-TextValue $ws.Range("E42") '  -3.23%  '
Set-TextValue $ws.Range("E43") '  -0.15%  '
Set-TextValue $ws.Range("E44") '  +8.02%  '
Set-TextValue $ws.Range("D45") '150.30'
Set-TextValue $ws.Range("E45") '  +1.76%  '
Set-TextValue $ws.Range("D46") '3.66'
Set-TextValue $ws.Range("E46") '  +1.45%  '
Set-TextValue $ws.Range("D47") '0.0539'
Set-TextValue $ws.Range("E47") '  +3.02%  '
Set-TextValue $ws.Range("D48") '20.34'
Set-TextValue $ws.Range("E48") '  +3.06%  '
Set-TextValue $ws.Range("D49") '0.594'
Set-TextValue $ws.Range("E49") '  +2.11%  '
Set-TextValue $ws.Range("D50") '0.0923'
Set-TextValue $ws.Range("E50") '  +1.82%  '
Set-TextValue $ws.Range("E51") '  +2.12%  '
